$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply the "datetime" number format used by A2/A3 to the whole of column A ---
# (mirrors the original author re-formatting column A, which is why A1's header cell
# also ends up carrying that style even though it holds text)
$ws.Columns("A").NumberFormat = "m/d/yy h:mm"

# --- Fix up A4: it previously used the short-date style (s=2); it now matches
#     the date+time style used by A2/A3/A5 (s=1) ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A4").Value = 42277

# --- Add the new bug report row (row 5) ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A5").Value = 42277.71875
$ws.Range("B5").Value = "QS text scraper features duplicate names"
$ws.Range("C5").Value = "Matt"
$ws.Range("E5").Value = "Ranking Scrape"

# --- Column width tweaks (closest reproducible values in this environment) ---
$ws.Columns("A").ColumnWidth = 15.666666666666666
$ws.Columns("B").ColumnWidth = 33.83333333333333
$ws.Columns("D").ColumnWidth = 11.5
$ws.Columns("E").ColumnWidth = 20.5

# --- Update selection to match the new active cell ---
$ws.Range("B6").Select() | Out-Null

# --- Printer / page setup ---
$ws.PageSetup.PaperSize = 9     # xlPaperA4... (A4=9 in Excel's XlPaperSize enum)
$ws.PageSetup.Orientation = 1   # xlPortrait

$excel.CutCopyMode = $false
